$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -9716
$ws.Range("N18").ClearContents()

$ws.Range("H62").Value = 2999.75
$ws.Range("J62").Value = 4800
$ws.Range("L62").Value = 4800
$ws.Range("N62").Value = -6048

$ws.Range("H64").Value = 4999.3335
$ws.Range("I64").Value = 4999.5
$ws.Range("K64").Value = 4999.5
$ws.Range("M64").Value = -4751.5

$ws.Range("H65").Value = 2999.75
$ws.Range("J65").Value = 4800
$ws.Range("L65").Value = 24000
$ws.Range("N65").Value = -30240

$ws.Range("H67").Value = 4999.3335
$ws.Range("I67").Value = 4999.5
$ws.Range("K67").Value = 4999.5
$ws.Range("M67").Value = -4141.5

$ws.Range("H69").Value = 8899

$ws.Range("H72").Value = 8899

$ws.Range("H76").Value = 20003188
$ws.Range("I76").Value = 33335784
$ws.Range("J76").Value = 4296.5
$ws.Range("K76").Value = 33335784
$ws.Range("L76").Value = 4296.5
$ws.Range("M76").Value = -33335469
$ws.Range("N76").Value = -4926.5

$ws.Range("H79").Value = 20003188
$ws.Range("I79").Value = 33335784
$ws.Range("J79").Value = 4296.5
$ws.Range("K79").Value = 33335784
$ws.Range("L79").Value = 4296.5
$ws.Range("M79").Value = -33334692
$ws.Range("N79").Value = -6480.5

$ws.Range("H100").Value = 2337.647
$ws.Range("I100").Value = 2189.4
$ws.Range("K100").Value = 2189.4
$ws.Range("M100").Value = -1648.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2233.3333
$ws.Range("I2").Value = 1900
$ws.Range("J2").Value = 2400
$ws.Range("K2").Value = 1900
$ws.Range("L2").Value = 2400
$ws.Range("M2").Value = -1787
$ws.Range("N2").Value = -2626

$ws.Range("H116").Value = 2233.3333
$ws.Range("I116").Value = 1900
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 1900
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 394
$ws.Range("N116").Value = -6988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2233.3333
$ws.Range("I3").Value = 1900
$ws.Range("J3").Value = 2400
$ws.Range("K3").Value = 1900
$ws.Range("L3").Value = 2400
$ws.Range("M3").Value = -1786
$ws.Range("N3").Value = -2628

$ws.Range("H86").Value = 3494.1333
$ws.Range("I86").Value = 3877.923
$ws.Range("K86").Value = 3877.923
$ws.Range("M86").Value = -2754.923

$ws.Range("H89").Value = 3494.1333
$ws.Range("I89").Value = 3877.923
$ws.Range("K89").Value = 19389.615
$ws.Range("M89").Value = -13773.615

$ws.Range("H94").Value = 390

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4079.5557
$ws.Range("I122").Value = 4451.357
$ws.Range("K122").Value = 13354.071
$ws.Range("M122").Value = -10904.071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 538.8
$ws.Range("I68").Value = 448.66666
$ws.Range("J68").Value = 674
$ws.Range("K68").Value = 1345.99998
$ws.Range("L68").Value = 2022
$ws.Range("M68").Value = -534.9999800000001
$ws.Range("N68").Value = -3644

$ws.Range("H71").Value = 538.8
$ws.Range("I71").Value = 448.66666
$ws.Range("J71").Value = 674
$ws.Range("K71").Value = 4037.99994
$ws.Range("L71").Value = 6066
$ws.Range("M71").Value = 18.0000600000003
$ws.Range("N71").Value = -14178

$ws.Range("H106").Value = 18498
$ws.Range("I106").Value = 9990
$ws.Range("J106").Value = 20625
$ws.Range("K106").Value = 29970
$ws.Range("L106").Value = 61875
$ws.Range("M106").Value = -29024
$ws.Range("N106").Value = -63767

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H113").Value = 2666.8333
$ws.Range("I113").Value = 2250.25
$ws.Range("K113").Value = 2250.25
$ws.Range("M113").Value = -80.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5001.0415
$ws.Range("I7").Value = 4796.8945
$ws.Range("K7").Value = 4796.8945
$ws.Range("M7").Value = -4684.8945

$ws.Range("H22").Value = 1257.1428
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -105

$ws.Range("H27").Value = 1257.1428
$ws.Range("I27").Value = 400
$ws.Range("K27").Value = 400
$ws.Range("M27").Value = -293

$ws.Range("H55").Value = 710.7
$ws.Range("I55").Value = 359.57144
$ws.Range("K55").Value = 359.57144
$ws.Range("M55").Value = -186.57144

$ws.Range("H61").Value = 1952.125
$ws.Range("I61").Value = 1945.7142
$ws.Range("J61").Value = 1997
$ws.Range("K61").Value = 1945.7142
$ws.Range("L61").Value = 1997
$ws.Range("M61").Value = -1743.7142
$ws.Range("N61").Value = -2401

$ws.Range("H93").Value = 999.25

$ws.Range("H113").Value = 1952.125
$ws.Range("I113").Value = 1945.7142
$ws.Range("J113").Value = 1997
$ws.Range("K113").Value = 1945.7142
$ws.Range("L113").Value = 1997
$ws.Range("M113").Value = 224.2858000000001
$ws.Range("N113").Value = -6337

$ws.Range("H122").Value = 3550.524
$ws.Range("I122").Value = 3247.9375
$ws.Range("J122").Value = 4518.8
$ws.Range("K122").Value = 9743.8125
$ws.Range("L122").Value = 13556.4
$ws.Range("M122").Value = -7293.8125
$ws.Range("N122").Value = -18456.4

$ws.Range("H126").Value = 5001.0415
$ws.Range("I126").Value = 4796.8945
$ws.Range("K126").Value = 14390.6835
$ws.Range("M126").Value = -11920.6835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5248.875
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5248.875
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10497.75
$ws.Range("N81").Value = -12619.75
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 5248.875
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5248.875
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 52488.75
$ws.Range("N84").Value = -63096.75
$ws.Range("M84").ClearContents()

$ws.Range("H113").Value = 212.77777
$ws.Range("I113").Value = 265.2
$ws.Range("J113").Value = 147.25
$ws.Range("K113").Value = 795.5999999999999
$ws.Range("L113").Value = 441.75
$ws.Range("M113").Value = 1374.4
$ws.Range("N113").Value = -4781.75
